$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must stay as literal
# text (matching the source inlineStr cells) - format as Text first so Excel
# doesn't auto-convert them to numeric values on assignment.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.767.59'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '2.116.38'
$ws.Range("E3").Value = '  +10.04%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '255.35'
$ws.Range("E5").Value = '  +2.56%  '
$ws.Range("D6").Value = '0.668'
$ws.Range("E6").Value = '  -3.88%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '46.41'
$ws.Range("E8").Value = '  +5.73%  '
$ws.Range("D9").Value = '62.24'
$ws.Range("E9").Value = '  +7.52%  '
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("E11").Value = '  -3.05%  '
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '14.59'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '2.421.32'
$ws.Range("E14").Value = '  +10.20%  '
$ws.Range("D15").Value = '0.851'
$ws.Range("E15").Value = '  +6.54%  '
$ws.Range("D16").Value = '2.110.97'
$ws.Range("E16").Value = '  +9.93%  '
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '36.767.24'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = '74.41'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").Value = '242.04'
$ws.Range("E22").Value = '  -4.02%  '
$ws.Range("D23").Value = '5.26'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  -8.02%  '
$ws.Range("D26").Value = '172.92'
$ws.Range("E26").Value = '  +3.02%  '
$ws.Range("D27").Value = '21.40'
$ws.Range("E27").Value = '  +13.75%  '
$ws.Range("D28").Value = '9.22'
$ws.Range("E28").Value = '  +4.10%  '
$ws.Range("D29").Value = '2.04'
$ws.Range("E29").Value = '  -9.69%  '
$ws.Range("E30").Value = '  -4.21%  '
$ws.Range("D31").Value = '22.67'
$ws.Range("E31").Value = '  +49.51%  '
$ws.Range("D32").Value = '4.57'
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '0.0959'
$ws.Range("E33").Value = '  +13.70%  '
$ws.Range("D34").Value = '0.0605'
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("D35").Value = '2.40'
$ws.Range("E35").Value = '  +19.39%  '
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -3.42%  '
$ws.Range("D39").Value = '0.914'
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("E40").Value = '  -7.94%  '
$ws.Range("E41").Value = '  +6.78%  '
$ws.Range("E42").Value = '  -2.86%  '
$ws.Range("D43").Value = '99.18'
$ws.Range("E43").Value = '  -5.28%  '
$ws.Range("D44").Value = '2.82'
$ws.Range("E44").Value = '  +17.04%  '
$ws.Range("D45").Value = '16.23'
$ws.Range("E45").Value = '  -5.90%  '
$ws.Range("D46").Value = '1.364.49'
$ws.Range("E46").Value = '  +1.28%  '
$ws.Range("D47").Value = '0.0838'
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("D48").Value = '2.31'
$ws.Range("E48").Value = '  -3.47%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.301.36'
$ws.Range("E49").Value = '  +9.81%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '6.89'
$ws.Range("E50").Value = '  +7.12%  '
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = '2.83'
$ws.Range("E51").Value = '  +1.47%  '
